$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Replace-ParagraphXml($findText, $paraInnerXml) {
    $r = $d.Content
    $found = $r.Find.Execute($findText)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    $r.Text = ""
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $paraInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark near the top (author email line)
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2) Rewrite the "Game Mechanics" bullet list entries
# ---------------------------------------------------------------------

# Up Arrow
Replace-ParagraphXml "If the Up Arrow is pressed, the Player jumps." (
  '<w:p w14:paraId="2028D2AF" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Up Arrow </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">press: </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player jumps.</w:t></w:r>' +
  '</w:p>'
)

# Right Arrow
Replace-ParagraphXml "If the Right Arrow is pressed, the Player moves to the right of the screen (forward)" (
  '<w:p w14:paraId="59273CE5" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t>Right Arrow pres</w:t></w:r>' +
    '<w:r><w:t>s:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player moves to the right of the screen (forward)</w:t></w:r>' +
  '</w:p>'
)

# Left Arrow
Replace-ParagraphXml "If the Left Arrow is pressed, the Player moves to the left of the screen (backward)" (
  '<w:p w14:paraId="5157F54E" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t>Left Arrow is press</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player moves to the left of the screen (backward)</w:t></w:r>' +
  '</w:p>'
)

# Down Arrow -> Space Button
Replace-ParagraphXml "If the Down Arrow is pressed --? Crouch???" (
  '<w:p w14:paraId="366AA28E" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Space Button press: </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player jumps</w:t></w:r>' +
  '</w:p>'
)

# W -> A Key Press
Replace-ParagraphXml "If the W is pressed on Level 2" (
  '<w:p w14:paraId="33EAE6AC" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">A </w:t></w:r>' +
    '<w:r><w:t>Key Press</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> on Level 2: </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player dodges raindrop by moving to the first segment of screen</w:t></w:r>' +
  '</w:p>'
)

# A -> S Key Press
Replace-ParagraphXml "If the A is pressed on Level 2" (
  '<w:p w14:paraId="55A01A17" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t>S Key Press</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> on Level 2: </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>the</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> Player dodges raindrop by moving to the second segment of screen</w:t></w:r>' +
  '</w:p>'
)

# ---------------------------------------------------------------------
# 3) Merge the "S" and "D" bullet paragraphs into a single "D Key Press" bullet
#    (keeps the "S" paragraph's identity, drops the old "D" paragraph entirely)
# ---------------------------------------------------------------------
$pS = $null
$pD = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "If the S is pressed on Level 2*") { $pS = $p }
    if ($t -like "If the D is pressed on Level 2*") { $pD = $p }
}
$mergeRange = $d.Range($pS.Range.Start, $pD.Range.End)
$mergeXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' +
  '<w:p w14:paraId="7A15DABD" w14:textId="77777777" w:rsidR="000347F3" w:rsidRDefault="000347F3" w:rsidP="000347F3">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:spacing w:line="256" w:lineRule="auto"/></w:pPr>' +
    '<w:r><w:t xml:space="preserve">D </w:t></w:r>' +
    '<w:r><w:t>Key Press</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> on Level 2: </w:t></w:r>' +
    '<w:r><w:t>the Player dodges raindrop by moving to the third segment of screen</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$mergeRange.InsertXML($mergeXml) | Out-Null

# ---------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark, now at the end of the merged "D" bullet
# ---------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("the Player dodges raindrop by moving to the third segment of screen")
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
